$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 keeps the same text value; re-setting it produces a fresh shared-string entry
$ws.Range("A1").Value = "HK_G_acc_G"

$values = @(
    71.586715867158674,
    71.586715867158674,
    71.586715867158674,
    71.217712177121768,
    71.586715867158674,
    71.586715867158674,
    73.431734317343171,
    73.800738007380076,
    74.538745387453872,
    73.800738007380076,
    72.32472324723247,
    73.062730627306266,
    74.538745387453872,
    74.538745387453872,
    74.538745387453872,
    74.169741697416967,
    74.169741697416967,
    74.538745387453872,
    73.062730627306266,
    73.431734317343171,
    73.800738007380076,
    70.110701107011081,
    70.110701107011081,
    69.741697416974162,
    73.062730627306266,
    72.693726937269375,
    73.800738007380076,
    74.907749077490777,
    74.538745387453872,
    74.169741697416967,
    69.741697416974162,
    70.110701107011081,
    71.586715867158674,
    70.479704797047972,
    70.479704797047972,
    75.645756457564573,
    69.003690036900366,
    69.372693726937271,
    70.479704797047972,
    74.169741697416967,
    74.169741697416967,
    74.907749077490777,
    74.169741697416967,
    70.479704797047972,
    71.217712177121768,
    70.479704797047972,
    74.538745387453872,
    70.848708487084863,
    73.431734317343171
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
